$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.682.68"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "3.080.71"
$ws.Range("E3").Value = "  +3.32%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'388.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.78%  "
$ws.Range("D6").Value = "'103.52"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").Value = "'0.545"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'0.587"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.82%  "
$ws.Range("D10").Value = "'37.04"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.93%  "
$ws.Range("E11").Value = "  +0.28%  "
$ws.Range("D12").Value = "'0.0864"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.69%  "
$ws.Range("D13").Value = "3.555.49"
$ws.Range("E13").Value = "  +3.31%  "
$ws.Range("D14").Value = "'18.76"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.87%  "
$ws.Range("D15").Value = "'7.81"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.39%  "
$ws.Range("D16").Value = "3.069.44"
$ws.Range("E16").Value = "  +2.34%  "
$ws.Range("D17").Value = "'0.979"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.58%  "
$ws.Range("D18").Value = "'10.74"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.72%  "
$ws.Range("D19").Value = "51.747.43"
$ws.Range("E19").Value = "  +0.57%  "
$ws.Range("D20").Value = "'3.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.25%  "
$ws.Range("E21").Value = "  -1.15%  "
$ws.Range("E22").Value = "  +0.62%  "
$ws.Range("D23").Value = "'70.34"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("D24").Value = "'268.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.63%  "
$ws.Range("D25").Value = "'3.17"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.55%  "
$ws.Range("D26").Value = "'8.21"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.05%  "
$ws.Range("D27").Value = "'27.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.59%  "
$ws.Range("D28").Value = "'7.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("D29").Value = "'0.171"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.10%  "
$ws.Range("E30").Value = "  +0.13%  "
$ws.Range("D31").Value = "'0.108"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.10%  "
$ws.Range("E32").Value = "  -0.43%  "
$ws.Range("D33").Value = "'34.87"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.74%  "
$ws.Range("E34").Value = "  +0.34%  "
$ws.Range("D35").Value = "'50.10"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.85%  "
$ws.Range("D36").Value = "'0.0450"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.31%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("E38").Value = "  +1.99%  "
$ws.Range("E39").Value = "  +9.19%  "
$ws.Range("E40").Value = "  +2.37%  "
$ws.Range("D41").Value = "'16.98"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.35%  "
$ws.Range("D42").Value = "'2.58"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.97%  "
$ws.Range("E43").Value = "  -0.73%  "
$ws.Range("D44").Value = "'125.85"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.35%  "
$ws.Range("E45").Value = "  -0.54%  "
$ws.Range("D46").Value = "'21.88"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.99%  "
$ws.Range("E47").Value = "  +3.04%  "
$ws.Range("D48").Value = "'2.46"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.48%  "
$ws.Range("D49").Value = "2.039.50"
$ws.Range("E49").Value = "  +0.56%  "
$ws.Range("D50").Value = "3.376.18"
$ws.Range("E50").Value = "  +2.98%  "
$ws.Range("D51").Value = "'0.208"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.34%  "
